# ---------------------------------------------------------------------------
# template, with marginal notes
#
# 1) The (only, empty) paragraph in the body gets the "Note" paragraph
#    style applied, and its single run (a lone space character) is removed
#    - the bookmark ("_GoBack") stays untouched.
# 2) Five new paragraph styles are added to the style sheet:
#      - NormalWeb           (Word's built-in "Normal (Web)")
#      - MarginNoteOutside   (custom)
#      - MarginNoteInside    (custom, based on MarginNoteOutside)
#      - MarginNoteLeft      (custom, based on MarginNoteOutside)
#      - MarginNoteRIght     (custom, based on MarginNoteOutside)
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1. Apply the "Note" style to the lone paragraph and drop its run -----

$p = $d.Paragraphs.Item(1)
$p.Range.Style = "Note"

# Remove the single space character that made up the paragraph's only run,
# leaving the bookmark (_GoBack) in place.
$d.Content.Find.Execute(" ", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# --- 2. Add the new styles -------------------------------------------------

# "Normal (Web)" is one of Word's own built-in style ids - asking for it by
# its display name lets Word mint it as the real "NormalWeb" style id.
$normalWeb = $d.Styles.Add("Normal (Web)", 1)
$normalWeb.BaseStyle = "Normal"
$normalWeb.Priority = 99
$normalWeb.UnhideWhenUsed = $true
$normalWeb.ParagraphFormat.SpaceAfter = 5
$normalWeb.ParagraphFormat.SpaceAfterAuto = $true
$normalWeb.ParagraphFormat.LineSpacingRule = 0
$normalWeb.Font.Name = "Times"
$normalWeb.Font.NameBi = "Times New Roman"
$normalWeb.Font.Size = 10
$normalWeb.Font.SizeBi = 10
$normalWeb.LanguageID = "en-GB"

$marginNoteOutside = $d.Styles.Add("MarginNoteOutside", 1)
$marginNoteOutside.BaseStyle = "Normal"
$marginNoteOutside.QuickStyle = $true
$marginNoteOutside.ParagraphFormat.LineSpacingRule = 0
$marginNoteOutside.Font.Size = 10

$marginNoteInside = $d.Styles.Add("MarginNoteInside", 1)
$marginNoteInside.BaseStyle = "MarginNoteOutside"
$marginNoteInside.QuickStyle = $true

$marginNoteLeft = $d.Styles.Add("MarginNoteLeft", 1)
$marginNoteLeft.BaseStyle = "MarginNoteOutside"
$marginNoteLeft.QuickStyle = $true

$marginNoteRIght = $d.Styles.Add("MarginNoteRIght", 1)
$marginNoteRIght.BaseStyle = "MarginNoteOutside"
$marginNoteRIght.QuickStyle = $true
